$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: Id
$ws.Range("A2").Value = 105475261

# P2: Lokalnamn
$ws.Range("P2").Value = "6320217,477, Sm"

# S2: Noggrannhet
$ws.Range("S2").Value = 5

# X2: Externid - remove content (clear cell)
$ws.Range("X2").ClearContents()

# AC2: Publik kommentar
$ws.Range("AC2").Value = ". Trädslag först skriver som -Ask."

# AW2: Rapportör
$ws.Range("AW2").Value = "Torbjörn Blixt"

# AX2: Observatörer
$ws.Range("AX2").Value = "Via Torbjörn Blixt"

# AY2: Projektnamn
$ws.Range("AY2").Value = "Skyddsvärda träd"
